# Append four new paragraphs (an empty spacer paragraph followed by three
# "meter/second" unit-conversion notes) to the end of the document body,
# after the existing closing paragraph ("...which I do not have.").
#
# Word keeps a "_GoBack" bookmark around the most recent edit location, so
# once we add new trailing content that bookmark has to end up wrapping the
# very last run we insert instead of sitting where it used to (right after
# "...which I do not have."). We replace the whole last paragraph (via
# Range.InsertXML, which substitutes the contents of the Range it is called
# on) with: the same paragraph text it already had, minus the old
# bookmark, followed by the new paragraphs, with the bookmark re-created at
# the very end.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$lastPara = $d.Paragraphs.Last
$target = $lastPara.Range

$newFragment = @"
<w:p $wNs w:rsidR="008464BA" w:rsidRPr="00B60BA1" w:rsidRDefault="008464BA"><w:r><w:t>If I could calculate the expected number of changes in one time unit, I could know for sure</w:t></w:r><w:r w:rsidR="008F7251"><w:t xml:space="preserve"> when I had a matrix with 1 ETU</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="00264282"><w:t xml:space="preserve"> However, this depends </w:t></w:r><w:r w:rsidR="0017455B"><w:t>upon the amino acid frequencies, which I do not have.</w:t></w:r></w:p><w:p $wNs/><w:p $wNs><w:r><w:t>1 meter/second</w:t></w:r></w:p><w:p $wNs><w:r><w:t>1 meter / 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>-4</w:t></w:r><w:r><w:t xml:space="preserve"> second (much faster!)</w:t></w:r></w:p><w:p $wNs><w:r><w:t>10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve"> meter / second</w:t></w:r><w:r><w:t xml:space="preserve"> (multiplied the number by 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>4</w:t></w:r><w:r><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$target.InsertXML($newFragment) | Out-Null
